$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "G2" = 2
    "H2" = 2.72
    "I2" = 4.45
    "J2" = 2.77
    "K2" = 1.78
    "L2" = 5
    "M2" = 1.15
    "N2" = 4.2
    "O2" = 1.62
    "P2" = 2.02
    "Q2" = 2.77
    "R2" = 1.34
    "S2" = 1.62
    "T2" = 2.02
    "U2" = 2.27
    "W2" = 4.6
    "X2" = 7.7
    "Y2" = 9.75
    "Z2" = 18.5
    "AA2" = 24
    "AB2" = 55
    "AC2" = 4.5
    "AH2" = 8.5
    "AJ2" = 16
    "AL2" = 60
    "AM2" = 80
    "AO2" = 11.25
    "AP2" = 28
    "AQ2" = 50
    "AR2" = 120
    "AS2" = 500
    "AT2" = 2
    "AU2" = 8.5
    "AW2" = 5.8
    "AX2" = 28
    "AY2" = 40

    "G3" = 1.95
    "I3" = 3.4
    "J3" = 2.63
    "X3" = 10
    "Y3" = 9
    "AE3" = 13
    "AH3" = 11
    "AJ3" = 12
    "AO3" = 11
    "AU3" = 7.5
    "AX3" = 19
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
